$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- TONAME / (unlabeled) / EMAIL data -------------------------------
# Rows 2-4 already existed (previously Elva/Lynn, Daphne/Boone); they now
# all describe "Shah" / "Imran" and keep pointing at the gmail address.
$ws.Range("C2").Value = "Shah"
$ws.Range("G2").Value = "Imran"
$ws.Range("H2").Value = "shahimrans64@gmail.com"

$ws.Range("C3").Value = "Shah"
$ws.Range("G3").Value = "Imran"
$ws.Range("H3").Value = "shahimrans64@gmail.com"

$ws.Range("C4").Value = "Shah"
$ws.Range("G4").Value = "Imran"
$ws.Range("H4").Value = "shahimrans64@gmail.com"

# Rows 5-7 are brand new, pointing at the new outlook address.
$ws.Range("C5").Value = "Shah"
$ws.Range("G5").Value = "Imran"
$ws.Range("H5").Value = "shah_imran_sust@outlook.com"

$ws.Range("C6").Value = "Shah"
$ws.Range("G6").Value = "Imran"
$ws.Range("H6").Value = "shah_imran_sust@outlook.com"

$ws.Range("C7").Value = "Shah"
$ws.Range("G7").Value = "Imran"
$ws.Range("H7").Value = "shah_imran_sust@outlook.com"

# ---- Hyperlinks --------------------------------------------------------
# Drop all existing hyperlinks and rebuild them in the order/rIds used by
# the target workbook: H2->rId1, H3->rId2, H7->rId3, H6->rId4, H4->rId5,
# H5->rId6.
$ws.Hyperlinks.Delete()

$null = $ws.Hyperlinks.Add($ws.Range("H2"), "mailto:shahimrans64@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("H3"), "mailto:shahimrans64@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("H7"), "mailto:shah_imran_sust@outlook.com")
$null = $ws.Hyperlinks.Add($ws.Range("H6"), "mailto:shah_imran_sust@outlook.com")
$null = $ws.Hyperlinks.Add($ws.Range("H4"), "mailto:shahimrans64@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("H5"), "mailto:shah_imran_sust@outlook.com")

# Hyperlinks.Add() re-styles the cell it is attached to; make sure every
# EMAIL cell keeps using the same "hyperlink" cell style (s=2 in the
# original workbook).
$hlStyle = $ws.Range("H2").Style
$ws.Range("H2").Style = $hlStyle
$ws.Range("H3").Style = $hlStyle
$ws.Range("H4").Style = $hlStyle
$ws.Range("H5").Style = $hlStyle
$ws.Range("H6").Style = $hlStyle
$ws.Range("H7").Style = $hlStyle

# ---- Selection -----------------------------------------------------------
$null = $ws.Range("H5").Select()

Write-Host "done"
